$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E are safe to assign directly as text (Excel will not
# auto-convert these values to numbers), but column D values frequently
# look like numbers, so force those cells to Text format before writing
# and restore the default (unstyled) cell style afterwards so the saved
# file has no stray numeric value / style index.

$priceCells = @{
    'D2' = '68.567.09'
    'D3' = '3.757.32'
    'D4' = '0.999'
    'D5' = '584.38'
    'D6' = '177.03'
    'D7' = '3.747.62'
    'D8' = '0.638'
    'D10' = '0.723'
    'D11' = '0.168'
    'D12' = '53.92'
    'D13' = '0.0000303'
    'D14' = '10.86'
    'D15' = '4.330.90'
    'D16' = '3.743.06'
    'D17' = '19.61'
    'D18' = '13.19'
    'D19' = '1.15'
    'D20' = '0.127'
    'D21' = '68.390.40'
    'D22' = '415.07'
    'D23' = '4.59'
    'D24' = '89.57'
    'D25' = '3.13'
    'D26' = '13.03'
    'D27' = '10.89'
    'D29' = '5.98'
    'D30' = '9.73'
    'D31' = '8.06'
    'D32' = '33.40'
    'D33' = '12.88'
    'D35' = '44.40'
    'D36' = '616.20'
    'D37' = '65.88'
    'D38' = '0.0₃0937'
    'D39' = '0.408'
    'D41' = '0.998'
    'D42' = '3.23'
    'D44' = '3.11'
    'D46' = '2.65'
    'D47' = '9.50'
    'D48' = '0.137'
    'D49' = '2.74'
    'D50' = '2.765.58'
    'D51' = '0.000268'
}

foreach ($ref in $priceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$ref]
    $cell.Style = "Normal"
}

# Plain text updates (coin names, links, percentage strings)
$textCells = @{
    'E2' = '  -6.42%  '
    'E3' = '  -5.33%  '
    'E4' = '  -0.07%  '
    'E5' = '  -4.15%  '
    'E6' = '  +5.10%  '
    'E7' = '  -5.37%  '
    'E8' = '  -6.30%  '
    'E9' = '  +0.01%  '
    'E10' = '  -5.85%  '
    'E11' = '  -9.21%  '
    'E12' = '  -3.71%  '
    'E13' = '  -9.71%  '
    'E14' = '  -3.59%  '
    'E15' = '  -5.87%  '
    'E16' = '  -5.73%  '
    'E17' = '  -4.00%  '
    'E18' = '  -6.99%  '
    'B19' = 'Polygon'
    'C19' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E19' = '  -6.31%  '
    'B20' = 'TRON'
    'C20' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'E20' = '  -2.73%  '
    'E21' = '  -6.54%  '
    'E22' = '  -5.49%  '
    'E23' = '  -4.87%  '
    'E24' = '  -6.32%  '
    'E25' = '  -7.21%  '
    'E26' = '  -8.21%  '
    'E27' = '  -1.56%  '
    'E28' = '  -5.74%  '
    'E29' = '  +0.35%  '
    'E30' = '  -6.99%  '
    'B31' = 'NEARProtocol'
    'C31' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E31' = '  +3.45%  '
    'B32' = 'EthereumClassic'
    'C32' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'E32' = '  -7.37%  '
    'E33' = '  -6.84%  '
    'E34' = '  -7.65%  '
    'B35' = 'InjectiveProtocol'
    'C35' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'E35' = '  -6.97%  '
    'B36' = 'Bittensor'
    'C36' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'E36' = '  -4.81%  '
    'B37' = 'OKB'
    'C37' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E37' = '  -6.40%  '
    'E38' = '  -11.00%  '
    'E39' = '  -4.96%  '
    'E40' = '  +0.23%  '
    'E41' = '  -0.28%  '
    'E42' = '  +4.15%  '
    'E43' = '  -5.74%  '
    'E44' = '  -8.48%  '
    'E45' = '  -7.09%  '
    'E46' = '  +3.42%  '
    'E47' = '  -9.44%  '
    'E48' = '  -7.73%  '
    'E49' = '  -14.40%  '
    'E50' = '  -1.37%  '
    'B51' = 'FLOKI'
    'C51' = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
    'E51' = '  -10.40%  '
}

foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}
